$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meters")

# Rename existing "Caresens-N *" entry to "Caresens N *" and apply the
# "Input" cell style (matches the new library's naming convention).
$ws.Range("G2").Value = "Caresens N *"
$ws.Range("G2").Style = "Input"

# New Caresens / misc meter rows appended below the existing column G data.
$ws.Range("G3").Value = "Caresens N POP"

$ws.Range("G4").Value = "Caresens N Mini"
$ws.Range("G4").Style = "Bad"

$ws.Range("G5").Value = "Caresens N Voice"

$ws.Range("G6").Value = "Caresens II"

$ws.Range("G7").Value = "Caresens POP"
$ws.Range("G7").Style = "Bad"

$ws.Range("G8").Value = "COOL"

$ws.Range("G9").Value = "alphacheck professional"

# Column G now holds longer strings (matching column D's width).
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Move the active selection to where the editor left off.
$null = $ws.Range("H8").Select()
